$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 169.41176
$ws.Range("I11").Value = 169.41176
$ws.Range("K11").Value = 169.41176
$ws.Range("M11").Value = -29.41175999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2897.6
$ws.Range("I62").Value = 2897.6
$ws.Range("K62").Value = 2897.6
$ws.Range("M62").Value = -2273.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 2897.6
$ws.Range("I65").Value = 2897.6
$ws.Range("K65").Value = 14488
$ws.Range("M65").Value = -11368

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H99").Value = 331.36365
$ws.Range("J99").Value = 1090.5
$ws.Range("L99").Value = 3271.5
$ws.Range("N99").Value = -6267.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1026.5454
$ws.Range("I100").Value = 699.125
$ws.Range("J100").Value = 1899.6666
$ws.Range("K100").Value = 699.125
$ws.Range("L100").Value = 1899.6666
$ws.Range("M100").Value = -158.125
$ws.Range("N100").Value = -2981.6666

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 1739.1765
$ws.Range("J103").Value = 1631
$ws.Range("L103").Value = 4893
$ws.Range("N103").Value = -6065

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 15319.5
$ws.Range("I111").Value = 15319.5
$ws.Range("K111").Value = 45958.5
$ws.Range("M111").Value = -42891.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H118").Value = 448
$ws.Range("I118").Value = 448
$ws.Range("K118").Value = 1344
$ws.Range("M118").Value = 313

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H134").Value = 67288.46000000001
$ws.Range("J134").Value = 67288.46000000001
$ws.Range("L134").Value = 67288.46000000001
$ws.Range("N134").Value = -77428.46000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2085.2974
$ws.Range("I61").Value = 1652.069
$ws.Range("J61").Value = 3655.75
$ws.Range("K61").Value = 1652.069
$ws.Range("L61").Value = 3655.75
$ws.Range("M61").Value = -1440.069
$ws.Range("N61").Value = -4079.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 3706.5264
$ws.Range("I102").Value = 2766.2144
$ws.Range("K102").Value = 2766.2144
$ws.Range("M102").Value = -1144.2144

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H133").Value = 56051.6
$ws.Range("I133").Value = 20000
$ws.Range("J133").Value = 65064.5
$ws.Range("K133").Value = 20000
$ws.Range("L133").Value = 65064.5
$ws.Range("N133").Value = -70124.5
$ws.Range("M133").Value = -17470

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2085.2974
$ws.Range("I136").Value = 1652.069
$ws.Range("J136").Value = 3655.75
$ws.Range("K136").Value = 4956.207
$ws.Range("L136").Value = 10967.25
$ws.Range("M136").Value = -2406.207
$ws.Range("N136").Value = -16067.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 111009
$ws.Range("J59").Value = 111009
$ws.Range("L59").Value = 111009
$ws.Range("N59").Value = -112703

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 8548004
$ws.Range("I107").Value = 25641816
$ws.Range("J107").Value = 1097.5
$ws.Range("K107").Value = 25641816
$ws.Range("L107").Value = 1097.5
$ws.Range("M107").Value = -25639896
$ws.Range("N107").Value = -4937.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2041.6
$ws.Range("I16").Value = 1905.6
$ws.Range("K16").Value = 1905.6
$ws.Range("M16").Value = -1618.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3292628.5
$ws.Range("I31").Value = 2072.8462
$ws.Range("J31").Value = 10422166
$ws.Range("K31").Value = 2072.8462
$ws.Range("L31").Value = 10422166
$ws.Range("M31").Value = -1777.8462
$ws.Range("N31").Value = -10422756

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3292628.5
$ws.Range("I34").Value = 2072.8462
$ws.Range("J34").Value = 10422166
$ws.Range("K34").Value = 2072.8462
$ws.Range("L34").Value = 10422166
$ws.Range("M34").Value = -1870.8462
$ws.Range("N34").Value = -10422570

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 90286.664
$ws.Range("J52").Value = 90286.664
$ws.Range("L52").Value = 90286.664
$ws.Range("N52").Value = -90874.664

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2938.5
$ws.Range("I105").Value = 2346.25
$ws.Range("K105").Value = 2346.25
$ws.Range("M105").Value = -599.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 2041.6
$ws.Range("I113").Value = 1905.6
$ws.Range("K113").Value = 1905.6
$ws.Range("M113").Value = 264.4000000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 15879996
$ws.Range("I132").Value = 6192.4165
$ws.Range("J132").Value = 37045068
$ws.Range("K132").Value = 18577.2495
$ws.Range("L132").Value = 111135204
$ws.Range("M132").Value = -16047.2495
$ws.Range("N132").Value = -111140264

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H138").Value = 69999
$ws.Range("J138").Value = 69999
$ws.Range("L138").Value = 69999
$ws.Range("N138").Value = -80279

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 131.73334
$ws.Range("J12").Value = 170
$ws.Range("L12").Value = 510
$ws.Range("N12").Value = -856

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 675
$ws.Range("J34").Value = 1375
$ws.Range("L34").Value = 4125
$ws.Range("N34").Value = -4293

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 765.6667
$ws.Range("J46").Value = 779
$ws.Range("L46").Value = 2337
$ws.Range("N46").Value = -2519

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H52").Value = 4895
$ws.Range("J52").Value = 4895
$ws.Range("L52").Value = 14685
$ws.Range("N52").Value = -15217

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H99").Value = 8125
$ws.Range("J99").Value = 12000
$ws.Range("L99").Value = 36000
$ws.Range("N99").Value = -40492

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H108").Value = 1189.1428
$ws.Range("I108").Value = 1189.1428
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 3567.4284
$ws.Range("L108").Value = 0
$ws.Range("M108").Value = -687.4284000000002
$ws.Range("N108").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9760.450000000001
$ws.Range("I70").Value = 4909
$ws.Range("K70").Value = 4909
$ws.Range("M70").Value = -4639

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 9760.450000000001
$ws.Range("I73").Value = 4909
$ws.Range("K73").Value = 4909
$ws.Range("M73").Value = -3973

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3172.3635
$ws.Range("I132").Value = 2299.6667
$ws.Range("J132").Value = 3499.625
$ws.Range("K132").Value = 6899.000100000001
$ws.Range("L132").Value = 10498.875
$ws.Range("M132").Value = -4369.000100000001
$ws.Range("N132").Value = -15558.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1845.6875
$ws.Range("I61").Value = 1748.6666
$ws.Range("J61").Value = 2369.6
$ws.Range("K61").Value = 1748.6666
$ws.Range("L61").Value = 2369.6
$ws.Range("M61").Value = -1546.6666
$ws.Range("N61").Value = -2773.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1632.6923
$ws.Range("I93").Value = 1710.0834
$ws.Range("K93").Value = 1710.0834
$ws.Range("M93").Value = -462.0834

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1845.6875
$ws.Range("I113").Value = 1748.6666
$ws.Range("J113").Value = 2369.6
$ws.Range("K113").Value = 1748.6666
$ws.Range("L113").Value = 2369.6
$ws.Range("M113").Value = 421.3334
$ws.Range("N113").Value = -6709.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3826.7666
$ws.Range("I132").Value = 3409.3635
$ws.Range("J132").Value = 4974.625
$ws.Range("K132").Value = 10228.0905
$ws.Range("L132").Value = 14923.875
$ws.Range("M132").Value = -7698.0905
$ws.Range("N132").Value = -19983.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 12999.5
$ws.Range("J26").Value = 18999
$ws.Range("L26").Value = 18999
$ws.Range("N26").Value = -19585

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1000000000
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 1000000000
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 2000000000
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -2000001082

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1375.4445
$ws.Range("I113").Value = 1396.4
$ws.Range("J113").Value = 1349.25
$ws.Range("K113").Value = 4189.200000000001
$ws.Range("L113").Value = 4047.75
$ws.Range("M113").Value = -2019.200000000001
$ws.Range("N113").Value = -8387.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6088.7334
$ws.Range("I132").Value = 5333.3076
$ws.Range("K132").Value = 15999.9228
$ws.Range("M132").Value = -13469.9228

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3960
$ws.Range("I136").Value = 3265.1538
$ws.Range("J136").Value = 6971
$ws.Range("K136").Value = 9795.4614
$ws.Range("L136").Value = 20913
$ws.Range("M136").Value = -7245.4614
$ws.Range("N136").Value = -26013

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H138").Value = 89370
$ws.Range("J138").Value = 89370
$ws.Range("L138").Value = 89370
$ws.Range("N138").Value = -99650
